# Refresh cryptocurrency price/volume figures (GitHub Actions data sync).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.709.63"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "'1.892.59"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'245.05"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.4920"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.2956"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "'0.06794"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'1.887.34"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "'90.81"
$ws.Range("E13").Value = "  +5.11%  "
$ws.Range("D14").Value = "'0.6792"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "'5.040"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").Value = "'30.681.76"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "'0.000007984"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'13.15"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "'2.131.63"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'4.823"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'189.43"
$ws.Range("E23").Value = "  +32.96%  "
$ws.Range("D24").Value = "'6.143"
$ws.Range("E24").Value = "  +4.39%  "
$ws.Range("D25").Value = "'9.377"
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").Value = "'155.73"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").Value = "'19.07"
$ws.Range("E27").Value = "  +12.38%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'1.394"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "'4.337"
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("D31").Value = "'0.09078"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").Value = "'0.7504"
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").Value = "'1.108"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "'2.774"
$ws.Range("E36").Value = "  +4.26%  "
$ws.Range("D37").Value = "'0.01838"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'2.687"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "'2.145"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "'0.9381"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "'0.4423"
$ws.Range("E41").Value = "  +4.49%  "
$ws.Range("D42").Value = "'105.39"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "'5.768"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "'7.593"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").Value = "'0.1341"
$ws.Range("E46").Value = "  +5.00%  "
$ws.Range("D47").Value = "'0.05869"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").Value = "'8.712"
$ws.Range("E48").Value = "  +5.23%  "
$ws.Range("E49").Value = "  +5.84%  "
$ws.Range("D50").Value = "'0.3933"
$ws.Range("E50").Value = "  +3.83%  "
$ws.Range("D51").Value = "'33.62"
$ws.Range("E51").Value = "  +2.42%  "
